$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.38
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 1.06
$ws.Range("K2").Value = 9.5
$ws.Range("L2").Value = 1.3
$ws.Range("M2").Value = 3.4
$ws.Range("N2").Value = 2
$ws.Range("P2").Value = 1.4
$ws.Range("Q2").Value = 2.75
$ws.Range("R2").Value = 1.73
$ws.Range("S2").Value = 2
$ws.Range("T2").Value = 8
$ws.Range("U2").Value = 11
$ws.Range("V2").Value = 9.5
$ws.Range("W2").Value = 23
$ws.Range("X2").Value = 19
$ws.Range("Y2").Value = 29
$ws.Range("Z2").Value = 9.5
$ws.Range("AA2").Value = 6
$ws.Range("AB2").Value = 13
$ws.Range("AC2").Value = 51
$ws.Range("AD2").Value = 201
$ws.Range("AE2").Value = 9.5
$ws.Range("AF2").Value = 15
$ws.Range("AH2").Value = 29
$ws.Range("AI2").Value = 23
$ws.Range("AJ2").Value = 34

# Row 13
$ws.Range("G13").Value = 1.33
$ws.Range("I13").Value = 8
$ws.Range("T13").Value = 8
$ws.Range("Y13").Value = 26
$ws.Range("AA13").Value = 11
$ws.Range("AG13").Value = 21

# Row 20
$ws.Range("G20").Value = 4.5
$ws.Range("H20").Value = 2.95
$ws.Range("I20").Value = 1.9
$ws.Range("L20").Value = 1.5
$ws.Range("M20").Value = 2.27
$ws.Range("N20").Value = 2.42
$ws.Range("O20").Value = 1.44
$ws.Range("P20").Value = 1.53
$ws.Range("Q20").Value = 2.2
$ws.Range("R20").Value = 2.12
$ws.Range("S20").Value = 1.57
$ws.Range("T20").Value = 9.5
$ws.Range("U20").Value = 24
$ws.Range("V20").Value = 15.5
$ws.Range("W20").Value = 90
$ws.Range("X20").Value = 55
$ws.Range("Y20").Value = 70
$ws.Range("Z20").Value = 6.2
$ws.Range("AA20").Value = 6
$ws.Range("AB20").Value = 19.5
$ws.Range("AE20").Value = 5
$ws.Range("AF20").Value = 7.5
$ws.Range("AG20").Value = 9
$ws.Range("AH20").Value = 16
$ws.Range("AI20").Value = 20
$ws.Range("AJ20").Value = 45

# Row 21
$ws.Range("G21").Value = 2.07
$ws.Range("H21").Value = 2.95
$ws.Range("I21").Value = 3.75
$ws.Range("L21").Value = 1.4
$ws.Range("M21").Value = 2.52
$ws.Range("N21").Value = 2.18
$ws.Range("O21").Value = 1.53
$ws.Range("P21").Value = 1.47
$ws.Range("Q21").Value = 2.35
$ws.Range("R21").Value = 1.88
$ws.Range("S21").Value = 1.72
$ws.Range("T21").Value = 6.1
$ws.Range("U21").Value = 9
$ws.Range("V21").Value = 8.75
$ws.Range("W21").Value = 19.5
$ws.Range("X21").Value = 19
$ws.Range("Y21").Value = 35
$ws.Range("Z21").Value = 7.1
$ws.Range("AA21").Value = 5.8
$ws.Range("AB21").Value = 16
$ws.Range("AC21").Value = 90
$ws.Range("AD21").Value = 800
$ws.Range("AE21").Value = 9
$ws.Range("AF21").Value = 19.5
$ws.Range("AG21").Value = 13
$ws.Range("AH21").Value = 60
$ws.Range("AI21").Value = 40
$ws.Range("AJ21").Value = 50

# Row 22
$ws.Range("J22").Value = 1.05
$ws.Range("L22").Value = 1.3

# Row 23
$ws.Range("G23").Value = 1.21
$ws.Range("H23").Value = 6
$ws.Range("I23").Value = 10.5
$ws.Range("R23").Value = 1.5
$ws.Range("S23").Value = 2.27
$ws.Range("T23").Value = 16.5
$ws.Range("U23").Value = 10.75
$ws.Range("V23").Value = 10.25
$ws.Range("W23").Value = 10.25
$ws.Range("X23").Value = 9.75
$ws.Range("Y23").Value = 17
$ws.Range("Z23").Value = 35
$ws.Range("AA23").Value = 16
$ws.Range("AB23").Value = 18
$ws.Range("AC23").Value = 45
$ws.Range("AD23").Value = 175
$ws.Range("AE23").Value = 65
$ws.Range("AF23").Value = 150
$ws.Range("AG23").Value = 37
$ws.Range("AH23").Value = 400
$ws.Range("AI23").Value = 110
$ws.Range("AJ23").Value = 55

# Row 26
$ws.Range("J26").Value = 1.04
$ws.Range("K26").Value = 8.5
$ws.Range("L26").Value = 1.2
$ws.Range("M26").Value = 4.05
$ws.Range("P26").Value = 1.34
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 2.27
$ws.Range("S26").Value = 1.57
$ws.Range("Z26").Value = 12
$ws.Range("AE26").Value = 23
$ws.Range("AH26").Value = 251
$ws.Range("AI26").Value = 126
$ws.Range("AJ26").Value = 101

# Row 33
$ws.Range("K33").Value = 10

# Row 39
$ws.Range("G39").Value = 2.27
$ws.Range("H39").Value = 3.35
$ws.Range("I39").Value = 3
$ws.Range("K39").Value = 7.9
$ws.Range("Q39").Value = 2.75
$ws.Range("S39").Value = 2.1
$ws.Range("T39").Value = 8.25
$ws.Range("W39").Value = 25
$ws.Range("X39").Value = 19.5
$ws.Range("Z39").Value = 7.9
$ws.Range("AA39").Value = 6.8
$ws.Range("AB39").Value = 14
$ws.Range("AC39").Value = 60
$ws.Range("AD39").Value = 450
$ws.Range("AE39").Value = 9.75
$ws.Range("AF39").Value = 17.5
$ws.Range("AG39").Value = 11.25
$ws.Range("AH39").Value = 40
$ws.Range("AI39").Value = 27

# Row 50
$ws.Range("G50").Value = 1.48
$ws.Range("I50").Value = 6
$ws.Range("N50").Value = 1.6
$ws.Range("O50").Value = 2.3
$ws.Range("R50").Value = 1.75
$ws.Range("S50").Value = 2
$ws.Range("T50").Value = 8.5
$ws.Range("U50").Value = 8
$ws.Range("W50").Value = 11
$ws.Range("Y50").Value = 21
$ws.Range("AB50").Value = 15
$ws.Range("AC50").Value = 41

# Row 66
$ws.Range("G66").Value = 3.8
$ws.Range("H66").Value = 3.7
$ws.Range("I66").Value = 1.83
$ws.Range("K66").Value = 8
$ws.Range("L66").Value = 1.25
$ws.Range("M66").Value = 3.6
$ws.Range("N66").Value = 1.75
$ws.Range("O66").Value = 1.98
$ws.Range("P66").Value = 1.36
$ws.Range("Q66").Value = 2.87
$ws.Range("R66").Value = 1.7
$ws.Range("S66").Value = 2.05
$ws.Range("T66").Value = 12.5
$ws.Range("U66").Value = 22
$ws.Range("V66").Value = 12.5
$ws.Range("W66").Value = 55
$ws.Range("X66").Value = 32
$ws.Range("Y66").Value = 35
$ws.Range("Z66").Value = 8
$ws.Range("AA66").Value = 7.2
$ws.Range("AB66").Value = 14
$ws.Range("AC66").Value = 60
$ws.Range("AD66").Value = 400
$ws.Range("AE66").Value = 7.9
$ws.Range("AF66").Value = 9
$ws.Range("AG66").Value = 8.25
$ws.Range("AH66").Value = 15
$ws.Range("AI66").Value = 14

